$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update predicted win probabilities (column D) and predicted winner flags (column C)
# for rows 2-49, reflecting the refreshed model/stats output.

$ws.Range("D2").Value = 0.5061865278691161
$ws.Range("D3").Value = 0.285254186739531
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0.4066244996036203
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.334839697973324
$ws.Range("D6").Value = 0.2851706269867384
$ws.Range("D7").Value = 0.6561830178676972
$ws.Range("D8").Value = 0.3132137953446109
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0.4707068870421452
$ws.Range("D10").Value = 0.3442388589541159
$ws.Range("D11").Value = 0.6594134689891749
$ws.Range("D12").Value = 0.5607308807961117
$ws.Range("D13").Value = 0.4309928922501575
$ws.Range("D14").Value = 0.754742643802558
$ws.Range("D15").Value = 0.676947843220277
$ws.Range("D16").Value = 0.2990349665790054
$ws.Range("D17").Value = 0.4344026008854314
$ws.Range("D18").Value = 0.3441866172563843
$ws.Range("D19").Value = 0.3340558046887701
$ws.Range("D20").Value = 0.3973583570344865
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0.4800284657417949
$ws.Range("D22").Value = 0.7186372509269419
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0.2963391473341788
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0.3697516838353506
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0.3558370735052496
$ws.Range("D26").Value = 0.4307276344560493
$ws.Range("D27").Value = 0.4141562170280294
$ws.Range("D28").Value = 0.3407227340678441
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 0.4153108717673581
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 0.4811360996897946
$ws.Range("D31").Value = 0.5720585528692237
$ws.Range("D32").Value = 0.5884682227738453
$ws.Range("D33").Value = 0.2944547764319218
$ws.Range("D34").Value = 0.2896306462791369
$ws.Range("D35").Value = 0.7347571437876892
$ws.Range("D36").Value = 0.326690411517926
$ws.Range("D37").Value = 0.6323375930328625
$ws.Range("D38").Value = 0.6256916101200154
$ws.Range("D39").Value = 0.3126862779048019
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 0.4796859488182399
$ws.Range("D41").Value = 0.5538527614336057
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = 0.5642785180876689
$ws.Range("D43").Value = 0.7315973348712407
$ws.Range("D44").Value = 0.6770118392654273
$ws.Range("D45").Value = 0.356595899209315
$ws.Range("D46").Value = 0.733078730766656
$ws.Range("D47").Value = 0.4872543923041986
$ws.Range("D48").Value = 0.4428331428981567
$ws.Range("D49").Value = 0.7602541417771941
